$d = $word.ActiveDocument

# Update the header date line
$d.Content.Find.Execute("2023-11-28 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-29 Wednesday", 2) | Out-Null

# Update each practice-table cell in place (row, col) -> new text
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "93÷7="
$t.Cell(1,2).Range.Text = "30÷5="
$t.Cell(1,3).Range.Text = "19÷4="
$t.Cell(1,4).Range.Text = "18÷2="
$t.Cell(1,5).Range.Text = "51÷7="
$t.Cell(5,1).Range.Text = "72÷5="
$t.Cell(5,2).Range.Text = "99÷7="
$t.Cell(5,3).Range.Text = "94÷9="
$t.Cell(5,4).Range.Text = "71÷7="
$t.Cell(5,5).Range.Text = "75÷7="
$t.Cell(9,1).Range.Text = "47÷6="
$t.Cell(9,2).Range.Text = "30÷8="
$t.Cell(9,3).Range.Text = "87÷8="
$t.Cell(9,4).Range.Text = "97÷5="
$t.Cell(9,5).Range.Text = "57÷3="
$t.Cell(13,1).Range.Text = "95÷5="
$t.Cell(13,2).Range.Text = "19÷6="
$t.Cell(13,3).Range.Text = "85÷8="
$t.Cell(13,4).Range.Text = "57÷2="
$t.Cell(13,5).Range.Text = "25÷9="
$t.Cell(17,1).Range.Text = "55÷5="
$t.Cell(17,2).Range.Text = "30÷6="
$t.Cell(17,3).Range.Text = "85÷6="
$t.Cell(17,4).Range.Text = "23÷3="
$t.Cell(17,5).Range.Text = "15÷6="
